$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.358081459999084
$ws.Range("B1").Value = 2.165539264678955
$ws.Range("C1").Value = 4.904549121856689
$ws.Range("D1").Value = 3.305068254470825
$ws.Range("E1").Value = 1.270143985748291
